$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1084
$ws.Range("I12").Value = 875.5
$ws.Range("K12").Value = 875.5
$ws.Range("M12").Value = -705.5
$ws.Range("H41").Value = 780.125
$ws.Range("J41").Value = 402
$ws.Range("L41").Value = 402
$ws.Range("N41").Value = -1282
$ws.Range("H57").Value = 38999
$ws.Range("I57").Value = 38999
$ws.Range("K57").Value = 116997
$ws.Range("M57").Value = -116498
$ws.Range("H64").Value = 5020
$ws.Range("J64").Value = 5020
$ws.Range("L64").Value = 5020
$ws.Range("N64").Value = -5516
$ws.Range("H67").Value = 5020
$ws.Range("J67").Value = 5020
$ws.Range("L67").Value = 5020
$ws.Range("N67").Value = -6736
$ws.Range("H88").Value = 1462.5
$ws.Range("I88").Value = 1950
$ws.Range("K88").Value = 1950
$ws.Range("M88").Value = -1544
$ws.Range("H91").Value = 1462.5
$ws.Range("I91").Value = 1950
$ws.Range("K91").Value = 1950
$ws.Range("M91").Value = -546
$ws.Range("H116").Value = 4121.857
$ws.Range("I116").Value = 3410.6
$ws.Range("J116").Value = 5900
$ws.Range("K116").Value = 3410.6
$ws.Range("L116").Value = 5900
$ws.Range("M116").Value = 31.40000000000009
$ws.Range("N116").Value = -12784
$ws.Range("H137").Value = 2365.2222
$ws.Range("I137").Value = 1439.4
$ws.Range("J137").Value = 2721.3076
$ws.Range("K137").Value = 4318.200000000001
$ws.Range("L137").Value = 8163.9228
$ws.Range("M137").Value = -1768.200000000001
$ws.Range("N137").Value = -13263.9228

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H61").Value = 2490.111
$ws.Range("I61").Value = 2128.2856
$ws.Range("J61").Value = 3756.5
$ws.Range("K61").Value = 2128.2856
$ws.Range("L61").Value = 3756.5
$ws.Range("M61").Value = -1916.2856
$ws.Range("N61").Value = -4180.5
$ws.Range("H102").Value = 15631738
$ws.Range("I102").Value = 31253476
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 31253476
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -31251854
$ws.Range("N102").Value = -13244
$ws.Range("H136").Value = 2490.111
$ws.Range("I136").Value = 2128.2856
$ws.Range("J136").Value = 3756.5
$ws.Range("K136").Value = 6384.8568
$ws.Range("L136").Value = 11269.5
$ws.Range("M136").Value = -3834.8568
$ws.Range("N136").Value = -16369.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 33338230
$ws.Range("I107").Value = 100001290
$ws.Range("J107").Value = 6699.2
$ws.Range("K107").Value = 100001290
$ws.Range("L107").Value = 6699.2
$ws.Range("M107").Value = -99999370
$ws.Range("N107").Value = -10539.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3087.6667
$ws.Range("J22").Value = 4443.6665
$ws.Range("L22").Value = 4443.6665
$ws.Range("N22").Value = -5143.6665
$ws.Range("H52").Value = 83186.664
$ws.Range("J52").Value = 83186.664
$ws.Range("L52").Value = 83186.664
$ws.Range("N52").Value = -83774.664
$ws.Range("H107").Value = 2560.2856
$ws.Range("I107").Value = 732.25
$ws.Range("K107").Value = 732.25
$ws.Range("M107").Value = 1187.75
$ws.Range("H122").Value = 1795.091
$ws.Range("I122").Value = 1775.6
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 5326.799999999999
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -2876.799999999999
$ws.Range("N122").Value = -10870

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 581.96
$ws.Range("I107").Value = 274.9
$ws.Range("J107").Value = 786.6667
$ws.Range("K107").Value = 824.6999999999999
$ws.Range("L107").Value = 2360.0001
$ws.Range("M107").Value = 1095.3
$ws.Range("N107").Value = -6200.0001
$ws.Range("H132").Value = 4546.6
$ws.Range("I132").Value = 4443.4
$ws.Range("K132").Value = 39990.6
$ws.Range("M132").Value = -37460.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 154.78378
$ws.Range("J2").Value = 630.7143
$ws.Range("L2").Value = 630.7143
$ws.Range("N2").Value = -856.7143
$ws.Range("H80").Value = 2183.4443
$ws.Range("I80").Value = 2142.1667
$ws.Range("J80").Value = 2266
$ws.Range("K80").Value = 2142.1667
$ws.Range("L80").Value = 2266
$ws.Range("M80").Value = -1144.1667
$ws.Range("N80").Value = -4262
$ws.Range("H83").Value = 2183.4443
$ws.Range("I83").Value = 2142.1667
$ws.Range("J83").Value = 2266
$ws.Range("K83").Value = 10710.8335
$ws.Range("L83").Value = 11330
$ws.Range("M83").Value = -5718.833500000001
$ws.Range("N83").Value = -21314
$ws.Range("H97").Value = 888.8
$ws.Range("I97").Value = 400
$ws.Range("K97").Value = 400
$ws.Range("M97").Value = 96
$ws.Range("H102").Value = 4035.2222
$ws.Range("I102").Value = 3536.1667
$ws.Range("J102").Value = 5033.3335
$ws.Range("K102").Value = 3536.1667
$ws.Range("L102").Value = 5033.3335
$ws.Range("M102").Value = -1914.1667
$ws.Range("N102").Value = -8277.333500000001
$ws.Range("H132").Value = 1688.15
$ws.Range("I132").Value = 1688.15
$ws.Range("K132").Value = 5064.450000000001
$ws.Range("M132").Value = -2534.450000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1224
$ws.Range("H40").Value = 2415.4614
$ws.Range("I40").Value = 2600.182
$ws.Range("K40").Value = 2600.182
$ws.Range("M40").Value = -2464.182
$ws.Range("H93").Value = 2774.5
$ws.Range("I93").Value = 800
$ws.Range("K93").Value = 800
$ws.Range("M93").Value = 448
$ws.Range("H122").Value = 2613.8462
$ws.Range("I122").Value = 2485.6667
$ws.Range("J122").Value = 2902.25
$ws.Range("K122").Value = 7457.000100000001
$ws.Range("L122").Value = 8706.75
$ws.Range("M122").Value = -5007.000100000001
$ws.Range("N122").Value = -13606.75
$ws.Range("H126").Value = 1000
$ws.Range("J126").Value = 1000
$ws.Range("L126").Value = 3000
$ws.Range("N126").Value = -7940
$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 35000
$ws.Range("N134").Value = -45140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30100
$ws.Range("J54").Value = 30100
$ws.Range("L54").Value = 30100
$ws.Range("N54").Value = -31140
$ws.Range("H55").Value = 1123.6
$ws.Range("I55").Value = 874
$ws.Range("J55").Value = 1290
$ws.Range("K55").Value = 874
$ws.Range("L55").Value = 1290
$ws.Range("M55").Value = -597
$ws.Range("N55").Value = -1844
